$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 3.4
$ws.Range("I2").Value = 3.25
$ws.Range("T2").Value = 10
$ws.Range("U2").Value = 12
$ws.Range("V2").Value = 9
$ws.Range("X2").Value = 15
$ws.Range("AD2").Value = 126
